$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds for row 5 (Coritiba vs Botafogo SP) ---
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 6.25
$ws.Range("J5").Value = 2.25
$ws.Range("L5").Value = 6.5
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
$ws.Range("Z5").Value = 11
$ws.Range("AD5").Value = 7
$ws.Range("AF5").Value = 67
$ws.Range("AH5").Value = 13
$ws.Range("AI5").Value = 29
$ws.Range("AK5").Value = 67
$ws.Range("AN5").Value = 3.4
$ws.Range("AO5").Value = 8.5
$ws.Range("AQ5").Value = 29
$ws.Range("AU5").Value = 9.5
$ws.Range("AV5").Value = 67
$ws.Range("AX5").Value = 7
$ws.Range("AY5").Value = 34
$ws.Range("BA5").Value = 126
$ws.Range("BB5").Value = 151

# --- Update odds for row 6 (Avai vs Ponte Preta) ---
$ws.Range("G6").Value = 1.55
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 7
$ws.Range("Y6").Value = 8.5
$ws.Range("AE6").Value = 21
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 34
$ws.Range("BA6").Value = 151

# --- Remove the El Salvador "Platense Municipal vs Cacahuatique" fixture (row 8) ---
# Deleting the entire row shifts every row below it up by one, which also
# accounts for the dimension shrinking from BD18 to BD17.
$ws.Rows(8).Delete()
